$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ensemble 2 (row 5): new "Initial Drawdown Aq" value, styled like a
#     highlighted figure (14pt black Courier New) ---
$ws.Range("C5").Value = -0.57617189999999996
$ws.Range("C5").Font.Color = 0
$ws.Range("C5").Font.Name = "Courier New"
$ws.Range("C5").Font.Size = 14
$ws.Rows.Item(5).RowHeight = 18

# --- Ensemble 3 (row 6): fill in the previously-blank value ---
$ws.Range("C6").Value = -9.8889999999999993

# --- Ensemble 4 (row 7) ---
$ws.Range("C7").Value = -2.137

# --- Ensemble 5 (row 8) ---
$ws.Range("C8").Value = -2.1480000000000001

# --- Ensemble 6 (row 9) ---
$ws.Range("C9").Value = -2.1480000000000001

# Update the active selection to reflect where the user ended up editing
$ws.Range("F8").Select() | Out-Null
